$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PT_OBJECTIVES = 'Fazer previsões e explicar acerca dos efeitos sobre as velocidades das reações em vista de: catálise, variação da temperatura, geometria de colisão e concentração dos reagentes. Sugerir leis de velocidade de reação de posse de dados apropriados acerca dos efeitos de concentração, estudar as reações de ordens zero, um e dois. Estudar a aplicação da equação de Arrhenius. Ampliar o entendimento do sentido espontâneo das reações químicas. Entender a função termodinâmica entropia e sua relação com as três leis da termodinâmica. Entender o sentido de energia livre como uma referência para o grau de afastamento entre o sistema e seu estado de equilíbrio. Utilizar a variação da energia livre padrão como uma ferramenta para calcular a constante de equilíbrio para determinado processo. Examinar o conceito de equilíbrio e definir a constante de equilíbrio. Aprender a escrever as expressões das constantes de equilíbrio para reações homogêneas e heterogêneas e interpretar o sentido do quociente reacional. Dar a expressão do produto de solubilidade para um sal e calcular os produtos de solubilidade a partir de solubilidades determinadas experimentalmente e prever se deve ou não ocorrer precipitação. Aplicar os critérios de precipitação. Determinar os efeitos do íon comum. Calcular o pH de uma solução aquosa para sistemas envolvendo ácido ou base fortes ou pelo menos um ácido fraco ou uma base fraca. Entender o funcionamento de soluções-tampão. Estudar as reações envolvendo a formação de íons complexos a partir de espécies muito solúveis e muito pouco solúveis. Realizar uma representação simples para um sistema químico de uma pilha eletrolítica ou “galvânica” observando as convenções-padrão para identificar os eletrodos. Aplicar os princípios de estequiometria aos processos eletroquímicos usando equações balanceadas de semi-reações e o valor da constante de Faraday. Calcular potenciais-padrão de pilhas usando valores tabelados. Energia livre de Gibbs e a equação de Nernst.'
$LARISSA = '5817330 - Larissa de Freitas'
$PEDRO = '1506103 - Pedro Carlos de Oliveira'
$PT_SHORT_SYLLABUS = 'Cinética Química, Termodinâmica e Equilíbrio, Equilíbrio Químico, Eletroquímica.'
$PT_FULL_SYLLABUS = 'Cinética Química: Energia de ativação e catalisadores. Efeito da concentração dos reagentes e temperatura sobre a velocidade de reações químicas. Equações diferenciais de velocidade de reação. Leis de velocidade para reações de 1ª e 2ª ordens. Equação de Arrhenius. Termodinâmica e equilíbrio: Primeira Lei da Termodinâmica. Calor de reação e energia interna. Entalpia e variação de Entalpia. Entropia. Segunda Lei da Termodinâmica. Entropia e Desordem. Terceira Lei da Termodinâmica. Variação de Entropia numa reação. Energia Livre Padrões e Critério de Espontaneidade. Relação entre variação de energia livre padrão e a constante de Equilíbrio. Equilíbrio Químico: Natureza do equilíbrio químico. Quociente de reação e constante de equilíbrio. Efeito da concentração de reagentes e temperatura sobre o equilíbrio (princípio de Le Chatelier). Equilíbrio ácido-base. Equilíbrios em solução aquosa: solubilidade e íons complexos. Eletroquímica: Semi-reações. Potenciais de eletrodo padrão. Potencial de células galvânicas. Relação entre variação de energia livre padrão e potencial de célula. Energia livre de Gibbs e a Equação de Nernst. Eletrólise e lei de Faraday.'
$METODO_TEXT = 'Serão realizadas duas provas escritas'
$CRITERIO_TEXT = 'NF = (P1 + P2*2)/3'
$NORMA_TEXT = 'Será realizada uma avaliação (P3) englobando toda a ementa. A média final será obtida conforme equação: MF= (NF+P3)/2.'
$BIBLIO_TEXT = '1) ATIKNS, P.; JONES, L. Princípios de Química, 5ªEdição, Ed. Bookman, 2012. 2) BRADY, J.; HUMISTON, G.E. Química Geral Volume II, 2ª Edição, Ed. LTC, 2005. 3) BRADY, J.E.; RUSSELL, J.W.; HOLUM, J.R. Química a matéria e suas transformações Volume II 3ª Edição, Ed. LTC, 2010. 4) BRADY, J.E.; SENESE, F., Química – A matéria e suas transformações Volume II, Ed. LTC, 5ªEdição, 2010. 5) BROWNN, S.L.; HOLME, T.A. Química geral aplicada à engenharia. São Paulo: Ed. Cengage Learning, 2010. 6) BROWN, T. L.; LEMAY, H.E.L.; Jr BURSTEN, B.E.; BURDGE, J.R. Química a ciência central. 9ª Edição, Ed. Pearson Prentice Hall, 2005. 7) CHANG, R.; GOLDSBY, K.A., Química, 11ª Edição, Ed. AMGH Editora Ltda, 2013 8) KOTZ, J.C.; TREICHEL, P.M.; WEAVER, G.C., Química Geral e Reações Químicas, Volume II, 6ª Edição, Ed. Cengage Learning, 2009. 9) KOTZ, J.C.; TREICHEL, P.M.; TOWNSEND, J. R.; TREICHEL, D.A., Química Geral e Reações Químicas, Volume II, 9ª Edição, Ed. Cengage Learning, 2016.'

# --- Row 10: fill in the Portuguese "Objetivos" text that was missing (previously
#     mistakenly duplicated the "Larissa de Freitas" string) ---
$ws.Range("B10").Value = $PT_OBJECTIVES
$ws.Range("C10").Value = $PT_OBJECTIVES

# --- Insert two new rows at 13 for the two "Docentes responsáveis" names,
#     pushing everything below down by two rows ---
$ws.Rows.Item(13).Resize(2).Insert()

# New row 13/14 should only contain B/C cells (no label in column A), matching
# the formatting used by column B (wrap, not bold) / column C (wrap, red)
$ws.Range("A13:A14").Clear()
$ws.Range("B13:B14").WrapText = $true
$ws.Range("B13:B14").Font.Bold = $false

$ws.Range("B13").Value = $LARISSA
$ws.Range("C13").Value = $LARISSA

$ws.Range("B14").Value = $PEDRO
$ws.Range("C14").Value = $PEDRO

# --- Row 15 (previously row 13): "Programa resumido:" short syllabus (PT) ---
$ws.Range("B15").Value = $PT_SHORT_SYLLABUS
$ws.Range("C15").Value = $PT_SHORT_SYLLABUS

# --- Row 16 (previously row 14): "Short syllabus:" (EN) already correct, no change ---

# --- Row 17 (previously row 15): "Programa:" full syllabus (PT) ---
$ws.Range("B17").Value = $PT_FULL_SYLLABUS
$ws.Range("C17").Value = $PT_FULL_SYLLABUS

# --- Row 18 (previously row 16): "Syllabus:" (EN) already correct, no change ---

# --- Row 20 (previously row 18): "Método:" ---
$ws.Range("B20").Value = $METODO_TEXT
$ws.Range("C20").Value = $METODO_TEXT

# --- Row 21 (previously row 19): "Critério:" ---
$ws.Range("B21").Value = $CRITERIO_TEXT
$ws.Range("C21").Value = $CRITERIO_TEXT

# --- Row 22 (previously row 20): "Norma de recuperação:" ---
$ws.Range("B22").Value = $NORMA_TEXT
$ws.Range("C22").Value = $NORMA_TEXT

# --- Row 23 (previously row 21): "Bibliografia:" ---
$ws.Range("B23").Value = $BIBLIO_TEXT
$ws.Range("C23").Value = $BIBLIO_TEXT
